{"js": "// Fix typo at p&p 613: \"ration\" -> \"ratio\"\n// (the sentence \"...leaves the triangle along its longest side equals to the\n// ration between the angle...\" should read \"...equals to the ratio between\n// the angle...\")\n\nconst results = context.document.body.search(\"the ration between the angle\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find text to fix ('the ration between the angle').\");\n}\n\nresults.items[0].insertText(\"the ratio between the angle\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix typo at p&p 613: \"ration\" -> \"ratio\"\n# (the sentence \"...leaves the triangle along its longest side equals to the\n# ration between the angle...\" should read \"...equals to the ratio between\n# the angle...\")\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"the ration between the angle\"\n$find.Replacement.Text = \"the ratio between the angle\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n"}
